# Insert a duplicate of the Betarraga "Primera/Segunda" record (rows 394:395,
# date 44217) as a new weekly entry right before the following record.
# This shifts the existing rows 396:516 down to 398:518 and grows the used
# range from A1:R516 to A1:R518, matching the new weekly price update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("394:395").Copy()
$ws.Rows("396:397").Insert()
